$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country name ordering (shared-string shuffle) and refreshed statistics
# for Pais worksheet, per commit "Update countries & provincias Spain"

$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 17:46"

$ws.Range("B6").Value = 40773
$ws.Range("C6").Value = 7227
$ws.Range("E6").Value = 40107
$ws.Range("G6").Value = 60
$ws.Range("H6").Value = 479

$ws.Range("B8").Value = 27558
$ws.Range("C8").Value = 2685
$ws.Range("E8").Value = 27021

$ws.Range("E13").Value = 5367
$ws.Range("G13").Value = 54
$ws.Range("H13").Value = 335

$ws.Range("B15").Value = 4306
$ws.Range("C15").Value = 724
$ws.Range("E15").Value = 4276

$ws.Range("A31").Value = "Luxemburgo"
$ws.Range("B31").Value = 875
$ws.Range("C31").Value = 77
$ws.Range("D31").Value = 6
$ws.Range("E31").Value = 861
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 8

$ws.Range("A32").Value = "Pakistan"
$ws.Range("B32").Value = 873
$ws.Range("C32").Value = 97
$ws.Range("D32").Value = 13
$ws.Range("E32").Value = 854
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 6

$ws.Range("A37").Value = "Grecia"
$ws.Range("B37").Value = 695
$ws.Range("C37").Value = 71
$ws.Range("D37").Value = 19
$ws.Range("E37").Value = 659
$ws.Range("F37").Value = 35
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 17

$ws.Range("A38").Value = "Polonia"
$ws.Range("B38").Value = 692
$ws.Range("C38").Value = 58
$ws.Range("D38").Value = 13
$ws.Range("E38").Value = 671
$ws.Range("F38").Value = 3
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 8

$ws.Range("A83").Value = "Republica de Chipre"
$ws.Range("B83").Value = 116
$ws.Range("C83").Value = 21
$ws.Range("D83").Value = 3
$ws.Range("E83").Value = 112
$ws.Range("F83").Value = 3
$ws.Range("H83").Value = 1

$ws.Range("A84").Value = "Jordania"
$ws.Range("B84").Value = 112
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 1
$ws.Range("E84").Value = 111
$ws.Range("F84").Value = 0
$ws.Range("H84").Value = 0

$ws.Range("A85").Value = "Moldavia"
$ws.Range("B85").Value = 109
$ws.Range("C85").Value = 15
$ws.Range("E85").Value = 106
$ws.Range("F85").Value = 3
$ws.Range("H85").Value = 1

$ws.Range("A86").Value = "Malta"
$ws.Range("B86").Value = 107
$ws.Range("C86").Value = 17
$ws.Range("E86").Value = 105
$ws.Range("F86").Value = 1
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0

$ws.Range("A87").Value = "Albania"
$ws.Range("B87").Value = 104
$ws.Range("C87").Value = 15
$ws.Range("D87").Value = 2
$ws.Range("E87").Value = 98
$ws.Range("F87").Value = 2
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 4

$ws.Range("A88").Value = "Nueva Zelanda"
$ws.Range("B88").Value = 102
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 102
$ws.Range("H88").Value = 0

$ws.Range("A89").Value = "Burkina Faso"
$ws.Range("B89").Value = 99
$ws.Range("C89").Value = 24
$ws.Range("D89").Value = 5
$ws.Range("E89").Value = 90
$ws.Range("F89").Value = 0
$ws.Range("H89").Value = 4

$ws.Range("A124").Value = "Paraguay"
$ws.Range("C124").Value = 0
$ws.Range("F124").Value = 1
$ws.Range("G124").Value = 0

$ws.Range("A125").Value = "Montenegro"
$ws.Range("C125").Value = 1
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 1

$ws.Range("A132").Value = "Togo"
$ws.Range("C132").Value = 2

$ws.Range("A133").Value = "Polinesia Francesa"
$ws.Range("C133").Value = 0

$ws.Range("A151").Value = "Haiti"
$ws.Range("C151").Value = 3

$ws.Range("A152").Value = "Surinam"
$ws.Range("C152").Value = 0

$ws.Range("A154").Value = "Bahamas"

$ws.Range("A155").Value = "Guinea"
$ws.Range("C155").Value = 2

$ws.Range("A156").Value = "Groenlandia"

$ws.Range("A157").Value = "Suazilandia"
$ws.Range("C157").Value = 0

$ws.Range("A159").Value = "El Salvador"

$ws.Range("A160").Value = "Namibia"

$ws.Range("A161").Value = "Republica de Africa Central"

$ws.Range("A162").Value = "Zambia"
$ws.Range("C162").Value = 0

$ws.Range("A163").Value = "Fiyi"
$ws.Range("C163").Value = 1

$ws.Range("A165").Value = "Cabo Verde"

$ws.Range("A166").Value = "Liberia"

$ws.Range("A167").Value = "Congo"

$ws.Range("A170").Value = "Benin"

$ws.Range("A171").Value = "Angola"

$ws.Range("A172").Value = "Niger"

$ws.Range("A173").Value = "Butan"

$ws.Range("A174").Value = "Nicaragua"

$ws.Range("A175").Value = "Santa Lucia"

$ws.Range("A176").Value = "Mauritania"

$ws.Range("A177").Value = "Sudan"
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 0
$ws.Range("H177").Value = 1

$ws.Range("A178").Value = "Nepal"
$ws.Range("D178").Value = 1
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

$ws.Range("A179").Value = "Gambia"
$ws.Range("C179").Value = 1
$ws.Range("G179").Value = 1

$ws.Range("A180").Value = "Islas Turcas y Caicos"
$ws.Range("C180").Value = 1

$ws.Range("A181").Value = "Santa Sede"

$ws.Range("A182").Value = "San Martin (Parte Holandesa)"

$ws.Range("A183").Value = "Papua Nueva Guinea"

$ws.Range("A184").Value = "Republica del Chad"

$ws.Range("A185").Value = "Timor Oriental"

$ws.Range("A186").Value = "Eritrea"

$ws.Range("A187").Value = "Republica de Yibuti"

$ws.Range("A188").Value = "Montserrat"
$ws.Range("C188").Value = 0

$ws.Range("A189").Value = "Uganda"

$ws.Range("A190").Value = "Mozambique"

$ws.Range("A191").Value = "Siria"

$ws.Range("A192").Value = "Granada"
$ws.Range("C192").Value = 0

$ws.Range("A193").Value = "Dominica"

$ws.Range("A194").Value = "Belice"
$ws.Range("C194").Value = 1

$ws.Range("A195").Value = "Somalia"

$ws.Range("A196").Value = "San Vicente y las Granadinas"

$ws.Range("A197").Value = "Antigua y Barbuda"
